$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'73343227"
$ws.Range("D2").Value = "'30677524"
